$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2999.8
$ws.Range("J70").Value = 3916.3333
$ws.Range("L70").Value = 11748.9999
$ws.Range("N70").Value = -12288.9999
$ws.Range("H73").Value = 2999.8
$ws.Range("J73").Value = 3916.3333
$ws.Range("L73").Value = 11748.9999
$ws.Range("N73").Value = -13620.9999
$ws.Range("H86").Value = 8980.833000000001
$ws.Range("I86").Value = 8224.5
$ws.Range("K86").Value = 8224.5
$ws.Range("M86").Value = -7101.5
$ws.Range("H89").Value = 8980.833000000001
$ws.Range("I89").Value = 8224.5
$ws.Range("K89").Value = 41122.5
$ws.Range("M89").Value = -35506.5
$ws.Range("H138").Value = 6381
$ws.Range("J138").Value = 7352.273
$ws.Range("L138").Value = 22056.819
$ws.Range("N138").Value = -32336.819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 3496.3635
$ws.Range("I28").Value = 3496.3635
$ws.Range("K28").Value = 3496.3635
$ws.Range("M28").Value = -3304.3635
$ws.Range("H61").Value = 3951.2273
$ws.Range("I61").Value = 2801.8708
$ws.Range("K61").Value = 2801.8708
$ws.Range("M61").Value = -2589.8708
$ws.Range("H63").Value = 3679.6875
$ws.Range("I63").Value = 2245.5454
$ws.Range("K63").Value = 2245.5454
$ws.Range("M63").Value = -1559.5454
$ws.Range("H66").Value = 3679.6875
$ws.Range("I66").Value = 2245.5454
$ws.Range("K66").Value = 11227.727
$ws.Range("M66").Value = -7795.726999999999
$ws.Range("H97").Value = 1340.5186
$ws.Range("I97").Value = 1450.7894
$ws.Range("J97").Value = 1078.625
$ws.Range("K97").Value = 1450.7894
$ws.Range("L97").Value = 1078.625
$ws.Range("M97").Value = -954.7893999999999
$ws.Range("N97").Value = -2070.625
$ws.Range("H99").Value = 3496.3635
$ws.Range("I99").Value = 3496.3635
$ws.Range("K99").Value = 3496.3635
$ws.Range("M99").Value = -501.3634999999999
$ws.Range("H136").Value = 3951.2273
$ws.Range("I136").Value = 2801.8708
$ws.Range("K136").Value = 8405.6124
$ws.Range("M136").Value = -5855.6124

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 58955.8
$ws.Range("J81").Value = 58955.8
$ws.Range("L81").Value = 58955.8
$ws.Range("N81").Value = -61077.8
$ws.Range("H84").Value = 58955.8
$ws.Range("J84").Value = 58955.8
$ws.Range("L84").Value = 176867.4
$ws.Range("N84").Value = -187475.4
$ws.Range("H94").Value = 66669890
$ws.Range("I94").Value = 2441.3333
$ws.Range("K94").Value = 2441.3333
$ws.Range("M94").Value = -1990.3333
$ws.Range("H134").Value = 3995.8462
$ws.Range("I134").Value = 2898.4075
$ws.Range("K134").Value = 8695.2225
$ws.Range("M134").Value = -6160.2225

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5188.1333
$ws.Range("I16").Value = 5712.1816
$ws.Range("K16").Value = 5712.1816
$ws.Range("M16").Value = -5425.1816
$ws.Range("H31").Value = 3714.3333
$ws.Range("I31").Value = 2670
$ws.Range("K31").Value = 2670
$ws.Range("M31").Value = -2375
$ws.Range("H34").Value = 3714.3333
$ws.Range("I34").Value = 2670
$ws.Range("K34").Value = 2670
$ws.Range("M34").Value = -2468
$ws.Range("H58").Value = 2649.7144
$ws.Range("J58").Value = 5599.75
$ws.Range("L58").Value = 5599.75
$ws.Range("N58").Value = -6005.75
$ws.Range("H86").Value = 9005.309999999999
$ws.Range("I86").Value = 12076.345
$ws.Range("K86").Value = 12076.345
$ws.Range("M86").Value = -10953.345
$ws.Range("H89").Value = 9005.309999999999
$ws.Range("I89").Value = 12076.345
$ws.Range("K89").Value = 60381.725
$ws.Range("M89").Value = -54765.725
$ws.Range("H113").Value = 5188.1333
$ws.Range("I113").Value = 5712.1816
$ws.Range("K113").Value = 5712.1816
$ws.Range("M113").Value = -3542.1816
$ws.Range("H134").Value = 3604.1428
$ws.Range("I134").Value = 2138.8572
$ws.Range("K134").Value = 6416.571599999999
$ws.Range("M134").Value = -3881.571599999999
$ws.Range("H136").Value = 2649.7144
$ws.Range("J136").Value = 5599.75
$ws.Range("L136").Value = 16799.25
$ws.Range("N136").Value = -21899.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5484.968
$ws.Range("J5").Value = 15261.8
$ws.Range("L5").Value = 45785.39999999999
$ws.Range("N5").Value = -46009.39999999999
$ws.Range("H29").Value = 6482573
$ws.Range("J29").Value = 649.75
$ws.Range("L29").Value = 1949.25
$ws.Range("N29").Value = -2503.25
$ws.Range("H44").Value = 3875
$ws.Range("I44").Value = 1250
$ws.Range("J44").Value = 6500
$ws.Range("K44").Value = 3750
$ws.Range("L44").Value = 19500
$ws.Range("M44").Value = -3352
$ws.Range("N44").Value = -20296
$ws.Range("H46").Value = 200229.8
$ws.Range("I46").Value = 200229.8
$ws.Range("K46").Value = 600689.3999999999
$ws.Range("M46").Value = -600598.3999999999
$ws.Range("H55").Value = 4981.909
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 4981.909
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 14945.727
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -15299.727
$ws.Range("H58").Value = 7749.875
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 8714.143
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 26142.429
$ws.Range("M58").Value = -2872
$ws.Range("N58").Value = -26398.429
$ws.Range("H132").Value = 3105.3333
$ws.Range("J132").Value = 3491.4614
$ws.Range("L132").Value = 31423.1526
$ws.Range("N132").Value = -36483.1526
$ws.Range("H133").Value = 5800.25
$ws.Range("I133").Value = 4486
$ws.Range("K133").Value = 13458
$ws.Range("M133").Value = -8398
$ws.Range("H134").Value = 10273.75
$ws.Range("J134").Value = 17142.857
$ws.Range("L134").Value = 51428.571
$ws.Range("N134").Value = -61568.571
$ws.Range("H135").Value = 5484.968
$ws.Range("J135").Value = 15261.8
$ws.Range("L135").Value = 137356.2
$ws.Range("N135").Value = -142426.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9727.9375
$ws.Range("I132").Value = 7345.421
$ws.Range("K132").Value = 22036.263
$ws.Range("M132").Value = -19506.263

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 494.6
$ws.Range("I16").Value = 494.6
$ws.Range("K16").Value = 494.6
$ws.Range("M16").Value = -324.6
$ws.Range("H55").Value = 1230.7916
$ws.Range("I55").Value = 1515.7059
$ws.Range("J55").Value = 538.8570999999999
$ws.Range("K55").Value = 1515.7059
$ws.Range("L55").Value = 538.8570999999999
$ws.Range("M55").Value = -1342.7059
$ws.Range("N55").Value = -884.8570999999999
$ws.Range("H61").Value = 2807.8064
$ws.Range("I61").Value = 2287.96
$ws.Range("J61").Value = 4973.8335
$ws.Range("K61").Value = 2287.96
$ws.Range("L61").Value = 4973.8335
$ws.Range("M61").Value = -2085.96
$ws.Range("N61").Value = -5377.8335
$ws.Range("H113").Value = 2807.8064
$ws.Range("I113").Value = 2287.96
$ws.Range("J113").Value = 4973.8335
$ws.Range("K113").Value = 2287.96
$ws.Range("L113").Value = 4973.8335
$ws.Range("M113").Value = -117.96
$ws.Range("N113").Value = -9313.833500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 37039260
$ws.Range("I81").Value = 1426.4286
$ws.Range("J81").Value = 166671660
$ws.Range("K81").Value = 2852.8572
$ws.Range("L81").Value = 333343320
$ws.Range("M81").Value = -1791.8572
$ws.Range("N81").Value = -333345442
$ws.Range("H84").Value = 37039260
$ws.Range("I84").Value = 1426.4286
$ws.Range("J84").Value = 166671660
$ws.Range("K84").Value = 14264.286
$ws.Range("L84").Value = 1666716600
$ws.Range("M84").Value = -8960.286
$ws.Range("N84").Value = -1666727208
$ws.Range("H122").Value = 3407.5925
$ws.Range("I122").Value = 2309.9048
$ws.Range("K122").Value = 6929.714399999999
$ws.Range("M122").Value = -4479.714399999999
